$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "25.846.18"
$ws.Range("E2").Value = "  +0.19%  "
$ws.Range("D3").Value = "1.734.30"
$ws.Range("E3").Value = "  -0.65%  "
$ws.Range("E4").Value = "  -0.07%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "231.48"
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.9999"
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.5149"
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.2781"
$ws.Range("E8").Value = "  +4.82%  "
$ws.Range("E9").Value = "  -2.55%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.06107"
$ws.Range("D11").Value = "1.753.22"
$ws.Range("E11").Value = "  -0.04%  "
$ws.Range("E12").Value = "  +1.32%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "15.24"
$ws.Range("E13").Value = "  -0.28%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "0.6434"
$ws.Range("E14").Value = "  +4.00%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "4.528"
$ws.Range("E15").Value = "  +1.53%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "76.83"
$ws.Range("E16").Value = "  -0.65%  "
$ws.Range("E17").Value = "  +0.00%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "1.0000"
$ws.Range("E18").Value = "  -0.10%  "
$ws.Range("D19").Value = "25.830.25"
$ws.Range("E19").Value = "  +0.06%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "11.51"
$ws.Range("E20").Value = "  -0.01%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "0.000006631"
$ws.Range("E21").Value = "  +0.75%  "
$ws.Range("D22").Value = "1.969.32"
$ws.Range("E22").Value = "  -0.28%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "4.149"
$ws.Range("E23").Value = "  +2.58%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "8.727"
$ws.Range("E24").Value = "  +6.13%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "5.121"
$ws.Range("E25").Value = "  +0.03%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "139.50"
$ws.Range("E26").Value = "  +2.75%  "
$ws.Range("E27").Value = "  +3.06%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "15.04"
$ws.Range("E28").Value = "  +0.32%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "1.795"
$ws.Range("E29").Value = "  +1.34%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "102.12"
$ws.Range("E30").Value = "  -0.41%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "0.08319"
$ws.Range("E31").Value = "  +1.76%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "3.698"
$ws.Range("E32").Value = "  +1.02%  "
$ws.Range("E33").Value = "  +1.76%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.04476"
$ws.Range("E34").Value = "  +2.45%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "2.616"
$ws.Range("E35").Value = "  -1.12%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.9812"
$ws.Range("E36").Value = "  -1.05%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.6138"
$ws.Range("E37").Value = "  +3.05%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "2.642"
$ws.Range("E38").Value = "  +0.80%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.01579"
$ws.Range("E39").Value = "  +1.65%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "1.946"
$ws.Range("E40").Value = "  +1.97%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.9992"
$ws.Range("E41").Value = "  -0.13%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "100.43"
$ws.Range("E42").Value = "  -0.87%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.3823"
$ws.Range("E43").Value = "  +0.42%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.7285"
$ws.Range("E44").Value = "  -2.21%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.05385"
$ws.Range("E46").Value = "  -1.74%  "
$ws.Range("B47").Value = "Algorand"
$ws.Range("C47").Value = "https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.1122"
$ws.Range("E47").Value = "  +2.80%  "
$ws.Range("B48").Value = "Aptos"
$ws.Range("C48").Value = "https://coinranking.com/coin/HGYj5JCv5+aptos-apt"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "6.238"
$ws.Range("E48").Value = "  +5.59%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "52.98"
$ws.Range("E49").Value = "  +1.33%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "30.03"
$ws.Range("E50").Value = "  +0.34%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "7.570"
$ws.Range("E51").Value = "  +2.43%  "
